$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking price
# strings (e.g. "58.47") are stored as text, matching the source data,
# then restore the default style so no stray formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "36.301.58"
$ws.Range("E2").Value = "  -2.96%  "

$ws.Range("D3").Value = "1.973.58"
$ws.Range("E3").Value = "  -3.38%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "245.09"
$ws.Range("E5").Value = "  -2.38%  "

$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  -4.06%  "

$ws.Range("D7").Value = "58.47"
$ws.Range("E7").Value = "  -9.78%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "0.372"
$ws.Range("E9").Value = "  -7.33%  "

$ws.Range("D10").Value = "56.38"
$ws.Range("E10").Value = "  -4.97%  "

$ws.Range("D11").Value = "0.0865"
$ws.Range("E11").Value = "  +9.47%  "

$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").Value = "0.851"
$ws.Range("E13").Value = "  -6.35%  "

$ws.Range("D14").Value = "21.87"
$ws.Range("E14").Value = "  -6.16%  "

$ws.Range("D15").Value = "2.271.53"
$ws.Range("E15").Value = "  -3.15%  "

$ws.Range("D16").Value = "13.64"
$ws.Range("E16").Value = "  -7.42%  "

$ws.Range("D17").Value = "5.42"
$ws.Range("E17").Value = "  -4.82%  "

$ws.Range("D18").Value = "1.981.89"
$ws.Range("E18").Value = "  -3.18%  "

$ws.Range("D19").Value = "36.243.62"
$ws.Range("E19").Value = "  -2.86%  "

$ws.Range("D20").Value = "0.0₃0901"
$ws.Range("E20").Value = "  +2.49%  "

$ws.Range("D21").Value = "70.16"
$ws.Range("E21").Value = "  -3.95%  "

$ws.Range("D22").Value = "5.23"
$ws.Range("E22").Value = "  -4.19%  "

$ws.Range("D23").Value = "233.98"
$ws.Range("E23").Value = "  -1.99%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").Value = "2.50"
$ws.Range("E25").Value = "  -3.56%  "

$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  -2.85%  "

$ws.Range("D27").Value = "9.74"
$ws.Range("E27").Value = "  -2.02%  "

$ws.Range("D28").Value = "165.41"
$ws.Range("E28").Value = "  +2.78%  "

$ws.Range("D29").Value = "19.83"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("D31").Value = "0.119"
$ws.Range("E31").Value = "  -2.21%  "

$ws.Range("D32").Value = "1.17"
$ws.Range("E32").Value = "  -0.79%  "

$ws.Range("D33").Value = "4.84"
$ws.Range("E33").Value = "  -5.57%  "

$ws.Range("D34").Value = "0.0645"
$ws.Range("E34").Value = "  +3.02%  "

$ws.Range("D35").Value = "4.39"
$ws.Range("E35").Value = "  -5.60%  "

$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").Value = "6.06"
$ws.Range("E37").Value = "  -5.23%  "

$ws.Range("E38").Value = "  -1.82%  "

$ws.Range("D39").Value = "2.19"
$ws.Range("E39").Value = "  -6.75%  "

$ws.Range("D40").Value = "2.90"
$ws.Range("E40").Value = "  -1.28%  "

$ws.Range("D41").Value = "0.0962"
$ws.Range("E41").Value = "  -5.26%  "

$ws.Range("E42").Value = "  -6.15%  "

$ws.Range("E43").Value = "  -5.26%  "

$ws.Range("E44").Value = "  -2.81%  "

$ws.Range("E45").Value = "  -6.84%  "

$ws.Range("D46").Value = "15.99"
$ws.Range("E46").Value = "  -7.31%  "

$ws.Range("D47").Value = "90.64"
$ws.Range("E47").Value = "  -4.53%  "

$ws.Range("D48").Value = "1.358.95"
$ws.Range("E48").Value = "  -2.42%  "

$ws.Range("D49").Value = "7.38"
$ws.Range("E49").Value = "  -5.48%  "

$ws.Range("E50").Value = "  -2.65%  "

$ws.Range("D51").Value = "45.08"
$ws.Range("E51").Value = "  -4.13%  "

$ws.Range("D2:D51").Style = "Normal"
